$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 226.8421
$ws.Range("I12").Value = 226.8421
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 226.8421
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -56.84209999999999
$ws.Range("N12").ClearContents()

$ws.Range("H28").Value = 508.8
$ws.Range("I28").Value = 262.5
$ws.Range("J28").Value = 1083.5
$ws.Range("K28").Value = 262.5
$ws.Range("L28").Value = 1083.5
$ws.Range("M28").Value = 222.5
$ws.Range("N28").Value = -2053.5

$ws.Range("H40").Value = 2749.2307
$ws.Range("I40").Value = 2464.4285
$ws.Range("J40").Value = 3081.5
$ws.Range("K40").Value = 2464.4285
$ws.Range("L40").Value = 3081.5
$ws.Range("M40").Value = -2289.4285
$ws.Range("N40").Value = -3431.5

$ws.Range("H86").Value = 49453.668
$ws.Range("I86").Value = 64434.188
$ws.Range("J86").Value = 1516
$ws.Range("K86").Value = 64434.188
$ws.Range("L86").Value = 1516
$ws.Range("M86").Value = -63311.188
$ws.Range("N86").Value = -3762

$ws.Range("H89").Value = 49453.668
$ws.Range("I89").Value = 64434.188
$ws.Range("J89").Value = 1516
$ws.Range("K89").Value = 322170.94
$ws.Range("L89").Value = 7580
$ws.Range("M89").Value = -316554.94
$ws.Range("N89").Value = -18812

$ws.Range("H132").Value = 5106840
$ws.Range("I132").Value = 5686992.5
$ws.Range("J132").Value = 1501.2
$ws.Range("K132").Value = 17060977.5
$ws.Range("L132").Value = 4503.6
$ws.Range("M132").Value = -17058447.5
$ws.Range("N132").Value = -9563.6

$ws.Range("H137").Value = 1296.0817
$ws.Range("I137").Value = 1158.5143
$ws.Range("J137").Value = 1640
$ws.Range("K137").Value = 3475.5429
$ws.Range("L137").Value = 4920
$ws.Range("M137").Value = -925.5429000000004
$ws.Range("N137").Value = -10020

$ws.Range("H138").Value = 4353.8184
$ws.Range("I138").Value = 4149
$ws.Range("J138").Value = 4599.6
$ws.Range("K138").Value = 12447
$ws.Range("L138").Value = 13798.8
$ws.Range("M138").Value = -7307
$ws.Range("N138").Value = -24078.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 9800
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 9800
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 9800
$ws.Range("N58").Value = -10660

$ws.Range("H74").Value = 1346.45
$ws.Range("I74").Value = 1312.0526
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1312.0526
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -438.0526
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 1346.45
$ws.Range("I77").Value = 1312.0526
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 6560.263
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -2192.263
$ws.Range("N77").Value = -18736

$ws.Range("H132").Value = 4254.2104
$ws.Range("I132").Value = 4500.3335
$ws.Range("J132").Value = 3650.0908
$ws.Range("K132").Value = 13501.0005
$ws.Range("L132").Value = 10950.2724
$ws.Range("M132").Value = -10971.0005
$ws.Range("N132").Value = -16010.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 291.66666
$ws.Range("I22").Value = 291.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 291.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -118.66666

$ws.Range("H99").Value = 1819.35
$ws.Range("I99").Value = 2195
$ws.Range("J99").Value = 1777.6111
$ws.Range("K99").Value = 2195
$ws.Range("L99").Value = 1777.6111
$ws.Range("M99").Value = -697
$ws.Range("N99").Value = -4773.6111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22123.328
$ws.Range("I31").Value = 1474.25
$ws.Range("J31").Value = 32896.76
$ws.Range("K31").Value = 1474.25
$ws.Range("L31").Value = 32896.76
$ws.Range("M31").Value = -1179.25
$ws.Range("N31").Value = -33486.76

$ws.Range("H34").Value = 22123.328
$ws.Range("I34").Value = 1474.25
$ws.Range("J34").Value = 32896.76
$ws.Range("K34").Value = 1474.25
$ws.Range("L34").Value = 32896.76
$ws.Range("M34").Value = -1272.25
$ws.Range("N34").Value = -33300.76

$ws.Range("H140").Value = 54800
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54800
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54800
$ws.Range("N140").Value = -65160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3500
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 10500
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -10724

$ws.Range("H5").Value = 1549.4783
$ws.Range("I5").Value = 1443.8
$ws.Range("J5").Value = 1630.7693
$ws.Range("K5").Value = 4331.4
$ws.Range("L5").Value = 4892.3079
$ws.Range("M5").Value = -4219.4
$ws.Range("N5").Value = -5116.3079

$ws.Range("H22").Value = 1149
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1149
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3447
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3785

$ws.Range("H27").Value = 1149
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1149
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3447
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3651

$ws.Range("H68").Value = 1624.6
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1624.6
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4873.799999999999
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6495.799999999999

$ws.Range("H71").Value = 1624.6
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1624.6
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14621.4
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22733.4

$ws.Range("H122").Value = 496.125
$ws.Range("I122").Value = 451.33334
$ws.Range("J122").Value = 523
$ws.Range("K122").Value = 4062.00006
$ws.Range("L122").Value = 4707
$ws.Range("M122").Value = -1612.00006
$ws.Range("N122").Value = -9607

$ws.Range("H135").Value = 1549.4783
$ws.Range("I135").Value = 1443.8
$ws.Range("J135").Value = 1630.7693
$ws.Range("K135").Value = 12994.2
$ws.Range("L135").Value = 14676.9237
$ws.Range("M135").Value = -10459.2
$ws.Range("N135").Value = -19746.9237

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10074.5
$ws.Range("I15").Value = 7200
$ws.Range("J15").Value = 10485.143
$ws.Range("K15").Value = 7200
$ws.Range("L15").Value = 10485.143
$ws.Range("M15").Value = -6912
$ws.Range("N15").Value = -11061.143

$ws.Range("H54").Value = 6923.1333
$ws.Range("I54").Value = 7070
$ws.Range("J54").Value = 6912.643
$ws.Range("K54").Value = 7070
$ws.Range("L54").Value = 6912.643
$ws.Range("M54").Value = -6550
$ws.Range("N54").Value = -7952.643
